# Refresh Slit3-Robo2 ligand-receptor edge metrics with updated TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.168788
$ws.Range("H2").Value = 3.506364
$ws.Range("I2").Value = 0.0139255825297802
$ws.Range("J2").Value = 0.01408364098536324
$ws.Range("M2").Value = 0.3615393333333333
$ws.Range("N2").Value = 1.084618
$ws.Range("O2").Value = 0.7649240942154193
$ws.Range("P2").Value = 0.7664955283791567
$ws.Range("Q2").Value = 0.4225628343279999
$ws.Range("R2").Value = 3.803065508951999
$ws.Range("S2").Value = 0.01065201360301419
$ws.Range("T2").Value = 0.01079504783857834

# Row 3
$ws.Range("G3").Value = 1.168788
$ws.Range("H3").Value = 3.506364
$ws.Range("I3").Value = 0.0139255825297802
$ws.Range("J3").Value = 0.01408364098536324
$ws.Range("O3").Value = 0.228925442648571
$ws.Range("P3").Value = 0.2293957393280025
$ws.Range("Q3").Value = 0.126464030388
$ws.Range("R3").Value = 1.138176273492
$ws.Range("S3").Value = 0.003187920144769139
$ws.Range("T3").Value = 0.003230727236267558

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 1.168788
$ws.Range("H4").Value = 3.506364
$ws.Range("I4").Value = 0.0139255825297802
$ws.Range("J4").Value = 0.01408364098536324
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.002907
$ws.Range("N4").Value = 0.005814
$ws.Range("O4").Value = 0.006150463136009796
$ws.Range("P4").Value = 0.004108732292840814
$ws.Range("Q4").Value = 0.003397666715999999
$ws.Range("R4").Value = 0.020386000296
$ws.Range("S4").Value = [double]"8.564878199687516E-05"
$ws.Range("T4").Value = [double]"5.786591051733835E-05"

# Row 5
$ws.Range("I5").Value = 0.9504675141158767
$ws.Range("J5").Value = 0.9612555315680539
$ws.Range("M5").Value = 0.3615393333333333
$ws.Range("N5").Value = 1.084618
$ws.Range("O5").Value = 0.7649240942154193
$ws.Range("P5").Value = 0.7664955283791567
$ws.Range("Q5").Value = 28.84132465141711
$ws.Range("R5").Value = 259.571921862754
$ws.Range("S5").Value = 0.7270355023162682
$ws.Range("T5").Value = 0.7367980665766426

# Row 6
$ws.Range("I6").Value = 0.9504675141158767
$ws.Range("J6").Value = 0.9612555315680539
$ws.Range("O6").Value = 0.228925442648571
$ws.Range("P6").Value = 0.2293957393280025
$ws.Range("S6").Value = 0.217586196392064
$ws.Range("T6").Value = 0.2205079233471858

# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("I7").Value = 0.9504675141158767
$ws.Range("J7").Value = 0.9612555315680539
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.002907
$ws.Range("N7").Value = 0.005814
$ws.Range("O7").Value = 0.006150463136009796
$ws.Range("P7").Value = 0.004108732292840814
$ws.Range("Q7").Value = 0.231902100357
$ws.Range("R7").Value = 1.391412602142
$ws.Range("S7").Value = 0.00584581540754457
$ws.Range("T7").Value = 0.003949541644225525

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.09388066666666665
$ws.Range("H8").Value = 0.2816419999999999
$ws.Range("I8").Value = 0.001118545854010694
$ws.Range("J8").Value = 0.001131241597962925
$ws.Range("M8").Value = 0.3615393333333333
$ws.Range("N8").Value = 1.084618
$ws.Range("O8").Value = 0.7649240942154193
$ws.Range("P8").Value = 0.7664955283791567
$ws.Range("Q8").Value = 0.03394155363955555
$ws.Range("R8").Value = 0.3054739827559999
$ws.Range("S8").Value = 0.0008556026742175431
$ws.Range("T8").Value = 0.0008670916263550734

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.09388066666666665
$ws.Range("H9").Value = 0.2816419999999999
$ws.Range("I9").Value = 0.001118545854010694
$ws.Range("J9").Value = 0.001131241597962925
$ws.Range("O9").Value = 0.228925442648571
$ws.Range("P9").Value = 0.2293957393280025
$ws.Range("Q9").Value = 0.010157982014
$ws.Range("R9").Value = 0.09142183812599998
$ws.Range("S9").Value = 0.0002560636047521221
$ws.Range("T9").Value = 0.0002595020027232961

# Row 10
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.09388066666666665
$ws.Range("H10").Value = 0.2816419999999999
$ws.Range("I10").Value = 0.001118545854010694
$ws.Range("J10").Value = 0.001131241597962925
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.002907
$ws.Range("N10").Value = 0.005814
$ws.Range("O10").Value = 0.006150463136009796
$ws.Range("P10").Value = 0.004108732292840814
$ws.Range("Q10").Value = 0.000272911098
$ws.Range("R10").Value = 0.001637466588
$ws.Range("S10").Value = [double]"6.879575041029371E-06"
$ws.Range("T10").Value = [double]"4.647968884555114E-06"

# Row 11
$ws.Range("G11").Value = 2.8258325
$ws.Range("H11").Value = 5.651664999999999
$ws.Range("I11").Value = 0.03366852131788238
$ws.Range("J11").Value = 0.0227004443433548
$ws.Range("M11").Value = 0.3615393333333333
$ws.Range("N11").Value = 1.084618
$ws.Range("O11").Value = 0.7649240942154193
$ws.Range("P11").Value = 0.7664955283791567
$ws.Range("Q11").Value = 1.021649598161666
$ws.Range("R11").Value = 6.129897588969999
$ws.Range("S11").Value = 0.02575386317265372
$ws.Range("T11").Value = 0.01739978908140138

# Row 12
$ws.Range("G12").Value = 2.8258325
$ws.Range("H12").Value = 5.651664999999999
$ws.Range("I12").Value = 0.03366852131788238
$ws.Range("J12").Value = 0.0227004443433548
$ws.Range("O12").Value = 0.228925442648571
$ws.Range("P12").Value = 0.2293957393280025
$ws.Range("Q12").Value = 0.3057579023324999
$ws.Range("R12").Value = 1.834547413995
$ws.Range("S12").Value = 0.007707581146019073
$ws.Range("T12").Value = 0.005207385213218048

# Row 13
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 2.8258325
$ws.Range("H13").Value = 5.651664999999999
$ws.Range("I13").Value = 0.03366852131788238
$ws.Range("J13").Value = 0.0227004443433548
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.002907
$ws.Range("N13").Value = 0.005814
$ws.Range("O13").Value = 0.006150463136009796
$ws.Range("P13").Value = 0.004108732292840814
$ws.Range("Q13").Value = 0.008214695077499999
$ws.Range("R13").Value = 0.03285878031
$ws.Range("S13").Value = 0.0002070769992095955
$ws.Range("T13").Value = [double]"9.327004873537746E-05"

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.06880966666666666
$ws.Range("H14").Value = 0.206429
$ws.Range("I14").Value = 0.0008198361824499672
$ws.Range("J14").Value = 0.0008291415052651543
$ws.Range("M14").Value = 0.3615393333333333
$ws.Range("N14").Value = 1.084618
$ws.Range("O14").Value = 0.7649240942154193
$ws.Range("P14").Value = 0.7664955283791567
$ws.Range("Q14").Value = 0.02487740101355555
$ws.Range("R14").Value = 0.2238966091219999
$ws.Range("S14").Value = 0.0006271124492655684
$ws.Range("T14").Value = 0.0006355332561793038

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.06880966666666666
$ws.Range("H15").Value = 0.206429
$ws.Range("I15").Value = 0.0008198361824499672
$ws.Range("J15").Value = 0.0008291415052651543
$ws.Range("O15").Value = 0.228925442648571
$ws.Range("P15").Value = 0.2293957393280025
$ws.Range("Q15").Value = 0.007445274742999998
$ws.Range("R15").Value = 0.06700747268699998
$ws.Range("S15").Value = 0.0001876813609666734
$ws.Range("T15").Value = 0.0001902015286078329

# Row 16
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.06880966666666666
$ws.Range("H16").Value = 0.206429
$ws.Range("I16").Value = 0.0008198361824499672
$ws.Range("J16").Value = 0.0008291415052651543
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 0.002907
$ws.Range("N16").Value = 0.005814
$ws.Range("O16").Value = 0.006150463136009796
$ws.Range("P16").Value = 0.004108732292840814
$ws.Range("Q16").Value = 0.000200029701
$ws.Range("R16").Value = 0.001200178206
$ws.Range("S16").Value = [double]"5.042372217725524E-06"
$ws.Range("T16").Value = [double]"3.406720478017581E-06"
